# Fix synergy card icon paths: prefix with "ui/assets/" to import
# topdown shooter assets under the new asset directory layout.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("__data")

$cells = @("I6", "I7", "I8", "I9")

foreach ($cellRef in $cells) {
    $range = $ws.Range($cellRef)
    $current = $range.Value2
    if ($current -notlike "ui/assets/*") {
        $range.Value2 = "ui/assets/" + $current
    }
}
